# Generate Report for Handoff
# Consolidates the e2e fixture set from 4 files (1 md + 2 png + config) down to
# 2 files (2 md + config), regenerating handoff status for the new files.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/0e8c021c604f158841c834592427a760259f71b6"
$zhBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/513bad8338c713e130d7261106a050f4d2d0dbc4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c1d1d5f8f76cbb93ca8cd0882cc6dc604ad72aa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$md1 = "7712ca7f-55cb-4ea6-8b16-8a3a3487afa2.md"
$md2 = "b07d73ed-43ea-4d82-90ed-f4f09c252c2b.md"
$cfg = ".localization-config"

$xlf1zh = "7712ca7f-55cb-4ea6-8b16-8a3a3487afa2.5ae0b5988e282534b84e59c077b189920b473614.zh-cn.xlf"
$xlf2zh = "b07d73ed-43ea-4d82-90ed-f4f09c252c2b.772bc25e666e30d138fa3394e3ca21f22c886e94.zh-cn.xlf"
$xlf1de = "7712ca7f-55cb-4ea6-8b16-8a3a3487afa2.5ae0b5988e282534b84e59c077b189920b473614.de-de.xlf"
$xlf2de = "b07d73ed-43ea-4d82-90ed-f4f09c252c2b.772bc25e666e30d138fa3394e3ca21f22c886e94.de-de.xlf"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$include         = "Include"
$ignored         = "Ignored"
$epoch           = "0001-01-01 00:00:00"
$dt1             = "2016-03-08 19:09:52"
$dt2             = "2016-03-08 19:09:59"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = $md1
$ws1.Range("B2").Value = $readyForHandoff
$ws1.Range("C2").Value = $readyForHandoff

$ws1.Range("A3").Value = $md2
$ws1.Range("B3").Value = $readyForHandoff
$ws1.Range("C3").Value = $readyForHandoff

$ws1.Range("A4").Value = $cfg
$ws1.Range("B4").Value = $notLocalized
$ws1.Range("C4").Value = $notLocalized

$ws1.Rows("5:5").Delete()

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$repoBase/e2e/$md1", "", "", $md1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$repoBase/e2e/$md2", "", "", $md2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$repoBase/$cfg", "", "", $cfg)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C2").Value = $xlf1zh
$ws2.Range("D2").Value = $dt1
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = $include

$ws2.Range("C3").Value = $xlf2zh
$ws2.Range("D3").Value = $dt1
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = $include
$ws2.Range("I3").ClearContents()

$ws2.Range("B4").Value = $notLocalized
$ws2.Range("C4").ClearFormats()
$ws2.Range("C4").ClearContents()
$ws2.Range("D4").Value = $epoch
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = $ignored
$ws2.Range("I4").ClearContents()

$ws2.Rows("5:5").Delete()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$repoBase/e2e/$md1", "", "", $md1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhBase/$xlf1zh", "", "", $xlf1zh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$repoBase/e2e/$md2", "", "", $md2)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhBase/$xlf2zh", "", "", $xlf2zh)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$repoBase/$cfg", "", "", $cfg)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("C2").Value = $xlf1de
$ws3.Range("D2").Value = $dt2
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = $include

$ws3.Range("C3").Value = $xlf2de
$ws3.Range("D3").Value = $dt2
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = $include
$ws3.Range("I3").ClearContents()

$ws3.Range("B4").Value = $notLocalized
$ws3.Range("C4").ClearFormats()
$ws3.Range("C4").ClearContents()
$ws3.Range("D4").Value = $epoch
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = $ignored
$ws3.Range("I4").ClearContents()

$ws3.Rows("5:5").Delete()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$repoBase/e2e/$md1", "", "", $md1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deBase/$xlf1de", "", "", $xlf1de)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$repoBase/e2e/$md2", "", "", $md2)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deBase/$xlf2de", "", "", $xlf2de)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$repoBase/$cfg", "", "", $cfg)
